$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: report time label changed (F24) ---
# Forced-text entry (leading apostrophe) to keep the quotePrefix/style
# identical to the original "16:31" cell, since "21:05" would otherwise
# be auto-interpreted as a time value.
$ws.Range("F24").Value = "'21:05"

# --- Row 1 height: recalculated slightly (58.5 -> 57.75) ---
$ws.Rows.Item(1).RowHeight = 57.75

# --- Updated supplier stock (J) and unit price (P) figures ---
# Column J = "Supplier Stock 1", column P = "Supplier Unit Price 1".
# Q (Adjusted Supplier Subtotal 1) recalculates automatically from P*O.
$ws.Range("J2").Value = 14550
$ws.Range("P2").Value = 0.019050000000000001

$ws.Range("J3").Value = 5760
$ws.Range("P3").Value = 0.03628

$ws.Range("J4").Value = 33732
$ws.Range("P4").Value = 0.01451

$ws.Range("J5").Value = 9007
$ws.Range("P5").Value = 0.029020000000000001

$ws.Range("J6").Value = 55905
$ws.Range("P6").Value = 0.068930000000000005

$ws.Range("J7").Value = 85016
$ws.Range("P7").Value = 0.10521

$ws.Range("J8").Value = 57670
$ws.Range("P8").Value = 0.27210000000000001

$ws.Range("J9").Value = 112788
$ws.Range("P9").Value = 0.34556999999999999

$ws.Range("J10").Value = 18449
$ws.Range("P10").Value = 0.67481000000000002

$ws.Range("J11").Value = 94693
$ws.Range("P11").Value = 0.24126

$ws.Range("J12").Value = 1804342
$ws.Range("P12").Value = 0.01179

$ws.Range("J13").Value = 68280
$ws.Range("P13").Value = 0.029020000000000001

$ws.Range("J14").Value = 83450
$ws.Range("P14").Value = 0.01179

$ws.Range("J15").Value = 131140
$ws.Range("P15").Value = 0.01179

$ws.Range("J16").Value = 88605
$ws.Range("P16").Value = 0.029020000000000001

$ws.Range("J17").Value = 857327
$ws.Range("P17").Value = 0.01451

$ws.Range("J18").Value = 16077
$ws.Range("P18").Value = 0.035369999999999999

$ws.Range("J19").Value = 14899
$ws.Range("P19").Value = 0.65395000000000003

# Row 20: only the stock quantity changed, price (3.47) stays the same.
$ws.Range("J20").Value = 806

$ws.Range("J21").Value = 27275
$ws.Range("P21").Value = 0.60768999999999995

# Row 22: only the stock quantity changed, price (1.23) stays the same.
$ws.Range("J22").Value = 478
